# Insert a new column before column B, so the existing PercActivations /
# PercSegmentAreas columns (currently B and C) shift to C and D, and the
# segment-name column (currently A) shifts to B. Column A will then be
# repurposed to hold a numeric index (0..18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift B:C -> C:D, keeping the labels/styles/data that were already there.
$ws.Columns("B:B").Insert()

# New header for the (now empty) column B: copy the header formatting
# (bold font + border) from the neighboring header cell, then set the text.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B1").Value = "segments"

# Segment names that used to live in column A, row by row, plus the new
# numeric index that will replace them in column A.
$names = @(
    "background",
    "back_bumper",
    "back_glass",
    "back_left_door",
    "back_left_light",
    "back_right_door",
    "back_right_light",
    "front_bumper",
    "front_glass",
    "front_left_door",
    "front_left_light",
    "front_right_door",
    "front_right_light",
    "hood",
    "left_mirror",
    "right_mirror",
    "tailgate",
    "trunk",
    "wheel"
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 1).Value = $i
}

# The column insert copied column A's bold/border style into the new
# column B data cells; the segment-name column should be unstyled (only
# the header row B1 is styled), matching the original layout of the
# PercActivations/PercSegmentAreas data columns.
$ws.Range("B2:B20").Style = "Normal"
